$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83.. down to 84..
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record
$ws.Cells.Item(83, 1).Value = 1
$ws.Cells.Item(83, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(83, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(83, 4).Value = 44615
$ws.Cells.Item(83, 5).Value = 15
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100102
$ws.Cells.Item(83, 8).Value = "Cítricos"
$ws.Cells.Item(83, 9).Value = 100102004
$ws.Cells.Item(83, 10).Value = "Mandarina"
$ws.Cells.Item(83, 11).Value = "Murcott"
$ws.Cells.Item(83, 12).Value = "Tercera"
$ws.Cells.Item(83, 13).Value = 250
$ws.Cells.Item(83, 14).Value = 14000
$ws.Cells.Item(83, 15).Value = 15000
$ws.Cells.Item(83, 16).Value = 14500
$ws.Cells.Item(83, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(83, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(83, 19).Value = 725
$ws.Cells.Item(83, 20).Value = 20
